$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.605.21"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.326.53"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "2.331.52"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.55%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "2.741.46"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "56.629.54"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "2.332.46"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "326.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.74%  "
$ws.Range("E28").Value = "  +9.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "0.0₃0738"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.914"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "141.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.379"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "277.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.560"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0218"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.69%  "
